# Fruta / hortaliza, semanal
# Rotate the weekly price-sheet rows: the data that used to sit in row 3
# moves up to row 2, row 4's data moves up to row 3, and row 2's original
# data wraps around to row 4. Columns A, B, C, E, F, G, H, O, R are identical
# across these rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 (was row 3)
$ws.Range("D2").Value = 44267
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = 1650
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("P2").Value = 275
$ws.Range("Q2").Value = 6

# New row 3 (was row 4)
$ws.Range("D3").Value = 44623
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1900
$ws.Range("N3").Value = "$/paquete"
$ws.Range("P3").Value = 1900
$ws.Range("Q3").Value = 1

# New row 4 (was row 2)
$ws.Range("D4").Value = 44377
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2800
$ws.Range("M4").Value = 2364
$ws.Range("N4").Value = "$/docena de matas"
$ws.Range("P4").Value = 394
$ws.Range("Q4").Value = 6
